# Apply Xhosa translations to the exit interview document.
$d = $word.ActiveDocument

function Replace-Text($old, $new, [bool]$wholeWord = $false) {
    $d.Content.Find.Execute(
        $old,      # FindText
        $true,     # MatchCase
        $wholeWord,# MatchWholeWord
        $false,    # MatchWildcards
        $false,    # MatchSoundsLike
        $false,    # MatchAllWordForms
        $true,     # Forward
        1,         # Wrap (wdFindContinue)
        $false,    # Format
        $new,      # ReplaceWith
        2          # Replace (wdReplaceAll)
    )
}

Replace-Text "Appendix 5: SWIFT Exit Interview Schedule" "ISihlomelo sesi-5: Ishedyuli yodliwano-ndlebe lokuphuma lwe-SWIFT"

Replace-Text "Hi! I just have three quick questions for you about your clinic visit today." "Molo! Ndinemibuzo ekhawulezileyo emithathu kuwe malunga notyelelo lwakho eklinikhi namhlanje."

Replace-Text "Are you a parent or a caregiver of a child? " "Ingaba ungumzali okanye umnonopheli womntwana?"

Replace-Text "Did you see a poster in the clinic today advertising a parenting programme?" "Ukhe wayibona ipowusta ekliniki namhlanje ibhengeza inkqubo yobuzali?"

Replace-Text "Document for interviewers to keep track of and tally responses:" "Yenza amaxwebhu odliwano-ndlebe ukuze ugcine umkhondo kunye nokuhlanganisa iimpendulo:"

Replace-Text "Site________________________      Clinic Name ____________________________   " "Indawo_________________________ Igama leklinikhi ______________________________________"

Replace-Text "RA ______________________________       Date_______________________ " "RA ___________________________________ Umhla___________________________________"

Replace-Text "Are you a parent or a caregiver of a child?" "Ingaba ungumzali okanye umnonopheli womntwana?"

Replace-Text "Did you see a poster in the clinic today advertising a parenting programme?" "Ukhe wayibona ipowusta ekliniki namhlanje ibhengeza inkqubo yobuzali?"

Replace-Text "Yes" "Ewe" $true
Replace-Text "No" "Hayi" $true
